# Trade #27 closed at 2026-02-17 13:19:07 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate stats now that trade #27 has closed
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.23
$wsSummary.Range("B4").Value = -0.78
$wsSummary.Range("B5").Value = -0.58
$wsSummary.Range("B6").Value = 27
$wsSummary.Range("B7").Value = 9
$wsSummary.Range("B9").Value = 33.33

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) stats
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.23
$wsStatus.Range("D4").Value = 27
$wsStatus.Range("E4").Value = -0.78
$wsStatus.Range("F4").Value = -0.77
$wsStatus.Range("G4").Value = 33.33

# ---------------------------------------------------------------------------
# All Trades + MarketMaking sheets: append new trade row (#27) as row 28
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 28

    $ws.Cells.Item($row, 1).Value = 27

    # Force the date to be stored as literal text (matching the existing
    # "Date" column cells) instead of being auto-converted to a date serial.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "13:19:01"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.110701
    $ws.Cells.Item($row, 7).Value = 0.19
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 71.6333
    $ws.Cells.Item($row, 10).Value = 0.08
    $ws.Cells.Item($row, 11).Value = 99.23
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
